# This presentation was edited purely by reordering six slides
# (positions 7-12 in the deck): the original sequence
#   7:"Librería ggplot2", 8:"Palabra clave", 9:"La vida es una batalla",
#   10:"Sintaxis" (código), 11:"Sintaxis" (explicación), 12:[picture]
# becomes
#   7:"Palabra clave", 8:"Librería ggplot2", 9:"Sintaxis" (código),
#   10:[picture], 11:"La vida es una batalla", 12:"Sintaxis" (explicación)
# No text/content inside any slide is changed - only slide order.

$p = $ppt.ActivePresentation

$p.Slides.Item(8).MoveTo(7)
$p.Slides.Item(10).MoveTo(9)
$p.Slides.Item(12).MoveTo(10)
